# Module Code Generator.xlsx - "Added IPV and 3 other resources"
#
# The generator inputs in row 2 are updated to the last-generated resource
# ("International Student and Scholar Services") and the generated code
# snippets for four resources are logged down column A (rows 9-19):
#   1) FHF_The Food Pantries Food Connect Map  (rows 9-11)
#   2) FHF_Capital Roots                       (rows 13-15)
#   3) Aca_International Student and Scholar Services (rows 17-19)
# A new, blank formatted row is appended at the bottom (row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the generator inputs in row 2 (last resource generated) ---
$ws.Range("A2").Value = "Aca_International Student and Scholar Services"
$ws.Range("B2").Value = "International Student and Scholar Services"
$ws.Range("C2").Value = "Acadata()"
$ws.Range("D2").Value = "International Student and Scholar Services"

$ws.Range("B2").Font.Name = "Times New Roman"
$ws.Range("B2").Font.Size = 12
$ws.Range("D2").Font.Name = "Times New Roman"
$ws.Range("D2").Font.Size = 12

# --- Log of generated module code: FHF_The Food Pantries Food Connect Map ---
$ws.Range("A9").Value = "mod_Accordion_ui('FHF_The Food Pantries Food Connect Map')"
$ws.Range("A10").Value = "mod_Accordion_server('FHF_The Food Pantries Food Connect Map', selector=selection, data=FHFdata(), title = c('The Food Pantries Food Connect Map'), Visible = T)"
$ws.Range("A11").Value = "mod_info_server('FHF_The Food Pantries Food Connect Map', selector = selection, data = FHFdata(), rownametitle = c('The Food Pantries Food Connect Map'), phone = F, website = T)"
$ws.Range("A12").Clear()

# --- Log of generated module code: FHF_Capital Roots ---
$ws.Range("A13").Value = "mod_Accordion_ui('FHF_Capital Roots')"
$ws.Range("A14").Value = "mod_Accordion_server('FHF_Capital Roots', selector=selection, data=FHFdata(), title = c('Capital Roots'), Visible = T)"
$ws.Range("A15").Value = "mod_info_server('FHF_Capital Roots', selector = selection, data = FHFdata(), rownametitle = c('Capital Roots'), phone = T, website = T)"

# --- Log of generated module code: Aca_International Student and Scholar Services ---
$ws.Range("A17").Value = "mod_Accordion_ui('Aca_International Student and Scholar Services')"
$ws.Range("A18").Value = "mod_Accordion_server('Aca_International Student and Scholar Services', selector=selection, data=Acadata(), title = c('International Student and Scholar Services'), Visible = T)"
$ws.Range("A19").Value = "mod_info_server('Aca_International Student and Scholar Services', selector = selection, data = Acadata(), rownametitle = c('International Student and Scholar Services'), phone = T, website = T)"

# --- New blank formatted row appended at the bottom ---
$row21 = $ws.Range("A21:I21")
$row21.Font.Name = "Times New Roman"
$row21.Font.Size = 12
$ws.Range("A21:B21").WrapText = $true
$ws.Range("D21:E21").WrapText = $true
$ws.Range("G21:H21").WrapText = $true

# --- Selection left on the new row, matching the editing session ---
$ws.Range("A21:XFD21").Select()
